$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87 (pushes existing rows 87..168 down to 88..169,
# inheriting formatting the way Excel's native row-insert does).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new weekly price record.
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 45167
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112031
$ws.Range("G87").Value = "Poroto verde"
$ws.Range("H87").Value = "Magnum"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 45
$ws.Range("K87").Value = 33000
$ws.Range("L87").Value = 33000
$ws.Range("M87").Value = 33000
$ws.Range("N87").Value = "`$/malla 25 kilos"
$ws.Range("O87").Value = "Perú"
$ws.Range("P87").Value = 1320
$ws.Range("Q87").Value = 25
$ws.Range("R87").Value = "Hortaliza"
